# Actualización automática 2025-10-31 17:30:09
#
# Applies the updated sales figures for CASTRO ALCIVAR EDA MARIA across the
# three report sheets: "VENTAS POR GRUPO", "VENTA MENSUAL" and
# "CUMPLIMIENTO MENSUAL". Downstream totals / percentage columns on the
# CUMPLIMIENTO MENSUAL sheet are recalculated to stay consistent with the
# new VENTA values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: VENTAS POR GRUPO (per-client sales by product group)
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("L31").Value = 743.41     # INTRIAGO ALVARADO BRENDA ALEJANDRA - PIEDRA SINTERIZADA
$wsGrupo.Range("M48").Value = 3834.35    # SALAZAR BALLADARES MARIA ANGELICA - PORCELANATO
$wsGrupo.Range("H52").Value = 635.96     # SISA GUANO CARLOS ALBERTO - INODOROS
$wsGrupo.Range("L52").Value = 4021.36    # SISA GUANO CARLOS ALBERTO - PIEDRA SINTERIZADA

# Row 60 holds the "<count> de 58" summary counters per column; bump the
# counters for the two columns that now have an additional non-zero entry.
$wsGrupo.Range("H60").Value = "3 de 58"
$wsGrupo.Range("L60").Value = "11 de 58"

# ---------------------------------------------------------------------
# Sheet: VENTA MENSUAL (per-client sales by month)
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F31").Value = 9123.92  # INTRIAGO ALVARADO BRENDA ALEJANDRA - octubre
$wsMensual.Range("F48").Value = 5784.88  # SALAZAR BALLADARES MARIA ANGELICA - octubre
$wsMensual.Range("F52").Value = 4657.32  # SISA GUANO CARLOS ALBERTO - octubre

# Row 60 is the "octubre" column total across all clients.
$wsMensual.Range("F60").Value = 91874.38

# ---------------------------------------------------------------------
# Sheet: CUMPLIMIENTO MENSUAL (compliance per product group)
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 6 - INODOROS
$wsCumpl.Range("D6").Value = 1288.16
$wsCumpl.Range("E6").Value = 1037.90694516821
$wsCumpl.Range("F6").Value = 0.5537931755041755

# Row 11 - PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 24699.04
$wsCumpl.Range("E11").Value = -5125.979750750303
$wsCumpl.Range("F11").Value = 1.261889540290298

# Row 12 - PORCELANATO
$wsCumpl.Range("D12").Value = 46111.11
$wsCumpl.Range("E12").Value = 2512.949999999997
$wsCumpl.Range("F12").Value = 0.9483187952630858

# Row 14 - TOTAL
$wsCumpl.Range("D14").Value = 97478.53999999999
$wsCumpl.Range("E14").Value = 2419.452841887858
$wsCumpl.Range("F14").Value = 0.9757807662290352
